$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C2").Value = 0.4009303330917504
$ws.Range("D2").Value = 0.1745796982821375
$ws.Range("E2").Value = 0.1699338955890113
$ws.Range("F2").Value = 1.610503378416844
$ws.Range("G2").Value = 0.9603159730419293
$ws.Range("H2").Value = 1.002447699743968
$ws.Range("J2").Value = 0.209174229140924
$ws.Range("K2").Value = 1.065926075380929
$ws.Range("L2").Value = 0.1537840619895121
$ws.Range("M2").Value = 0.416106485863672
$ws.Range("N2").Value = 1.500206560469479
$ws.Range("O2").Value = 3.971363781301136

$ws.Range("C3").Value = 0.3984785810648361
$ws.Range("D3").Value = 0.1728796381864299
$ws.Range("E3").Value = 0.170180072062438
$ws.Range("F3").Value = 1.616263491917032
$ws.Range("G3").Value = 0.9642966084793443
$ws.Range("H3").Value = 1.008607703329005
$ws.Range("J3").Value = 0.2104666141971023
$ws.Range("K3").Value = 0.9695258052293809
$ws.Range("L3").Value = 0.1544676840858408
$ws.Range("M3").Value = 0.3965333436124752
$ws.Range("N3").Value = 1.506027467533798
$ws.Range("O3").Value = 3.992397509057724

$ws.Range("C4").Value = 0.3971472061512173
$ws.Range("D4").Value = 0.171892555591505
$ws.Range("E4").Value = 0.1703889466852733
$ws.Range("F4").Value = 1.620581495991757
$ws.Range("G4").Value = 0.967296138267514
$ws.Range("H4").Value = 1.012796528086938
$ws.Range("J4").Value = 0.2113307714535679
$ws.Range("K4").Value = 0.9103731807836368
$ws.Range("L4").Value = 0.1549223329891518
$ws.Range("M4").Value = 0.3846158950277072
$ws.Range("N4").Value = 1.510092423985505
$ws.Range("O4").Value = 4.007330127077736

$ws.Range("C5").Value = 0.3966485302715199
$ws.Range("D5").Value = 0.1715046582659667
$ws.Range("E5").Value = 0.1704886129963405
$ws.Range("F5").Value = 1.622537713083645
$ws.Range("G5").Value = 0.9686581148891662
$ws.Range("H5").Value = 1.014605827591055
$ws.Range("J5").Value = 0.2117007039470842
$ws.Range("K5").Value = 0.886279295336351
$ws.Range("L5").Value = 0.1551164013750252
$ws.Range("M5").Value = 0.379785142466865
$ws.Range("N5").Value = 1.511872636260513
$ws.Range("O5").Value = 4.013922827880435

$ws.Range("C6").Value = 0.3965683785194329
$ws.Range("D6").Value = 0.1714411168964602
$ws.Range("E6").Value = 0.17050604217944
$ws.Range("F6").Value = 1.622874419397284
$ws.Range("G6").Value = 0.9688927037303188
$ws.Range("H6").Value = 1.014912442805489
$ws.Range("J6").Value = 0.211763205541363
$ws.Range("K6").Value = 0.8822792604079552
$ws.Range("L6").Value = 0.1551491580492073
$ws.Range("M6").Value = 0.3789845646276646
$ws.Range("N6").Value = 1.512175718923046
$ws.Range("O6").Value = 4.015048198989689

$ws.Range("C7").Value = 0.3971403030537317
$ws.Range("D7").Value = 0.1718872660839779
$ws.Range("E7").Value = 0.1703902318699431
$ws.Range("F7").Value = 1.620607082029338
$ws.Range("G7").Value = 0.9673139409705769
$ws.Range("H7").Value = 1.012820514530858
$ws.Range("J7").Value = 0.2113356884676421
$ws.Range("K7").Value = 0.910048194053104
$ws.Range("L7").Value = 0.1549249146264895
$ws.Range("M7").Value = 0.3845506411718276
$ws.Range("N7").Value = 1.510115931294969
$ws.Range("O7").Value = 4.007416983363228

$ws.Range("C8").Value = 0.4000489196735089
$ws.Range("D8").Value = 0.1739817801922712
$ws.Range("E8").Value = 0.1700068171651754
$ws.Range("F8").Value = 1.61232739341969
$ws.Range("G8").Value = 0.961573210446673
$ws.Range("H8").Value = 1.004487339829325
$ws.Range("J8").Value = 0.2096051976362965
$ws.Range("K8").Value = 1.032680552541109
$ws.Range("L8").Value = 0.1540125422344865
$ws.Range("M8").Value = 0.4093370463137731
$ws.Range("N8").Value = 1.502111883090066
$ws.Range("O8").Value = 3.97819749065593

$ws.Range("C9").Value = 0.4071291033789493
$ws.Range("D9").Value = 0.1785365826701906
$ws.Range("E9").Value = 0.1697115572677603
$ws.Range("F9").Value = 1.602285068398331
$ws.Range("G9").Value = 0.9547241355235485
$ws.Range("H9").Value = 0.9913682619111626
$ws.Range("J9").Value = 0.2067711757217836
$ws.Range("K9").Value = 1.273390450161855
$ws.Range("L9").Value = 0.1524994876245795
$ws.Range("M9").Value = 0.4587248717699737
$ws.Range("N9").Value = 1.490299839073856
$ws.Range("O9").Value = 3.936904876889116

$ws.Range("C10").Value = 0.4131652486892676
$ws.Range("D10").Value = 0.1821523669209881
$ws.Range("E10").Value = 0.1697712970873297
$ws.Range("F10").Value = 1.598678189074903
$ws.Range("G10").Value = 0.9523830824091988
$ws.Range("H10").Value = 0.9836895480214594
$ws.Range("J10").Value = 0.205028848281561
$ws.Range("K10").Value = 1.45030296549885
$ws.Range("L10").Value = 0.1515550524176454
$ws.Range("M10").Value = 0.4954690453362716
$ws.Range("N10").Value = 1.483975223169551
$ws.Range("O10").Value = 3.916323648605186

$ws.Range("C11").Value = 0.4160914717812716
$ws.Range("D11").Value = 0.1838550787617237
$ws.Range("E11").Value = 0.1698582013612793
$ws.Range("F11").Value = 1.597855302112819
$ws.Range("G11").Value = 0.9519032142055295
$ws.Range("H11").Value = 0.980620955855187
$ws.Range("J11").Value = 0.2043097514864662
$ws.Range("K11").Value = 1.530783641509231
$ws.Range("L11").Value = 0.1511614771022742
$ws.Range("M11").Value = 0.5122810717597019
$ws.Range("N11").Value = 1.481606259518983
$ws.Range("O11").Value = 3.909079161049249

$ws.Range("C12").Value = 0.4172253904036154
$ws.Range("D12").Value = 0.184508105451016
$ws.Range("E12").Value = 0.1698996671546027
$ws.Range("F12").Value = 1.597661206403515
$ws.Range("G12").Value = 0.95180568219007
$ws.Range("H12").Value = 0.9795199258311129
$ws.Range("J12").Value = 0.204047996723844
$ws.Range("K12").Value = 1.561258255345763
$ws.Range("L12").Value = 0.151017606144638
$ws.Range("M12").Value = 0.5186609069903625
$ws.Range("N12").Value = 1.480782031359198
$ws.Range("O12").Value = 3.90664037814463

$ws.Range("C13").Value = 0.4169800341211669
$ws.Range("D13").Value = 0.1843670988971695
$ws.Range("E13").Value = 0.1698903566840642
$ws.Range("F13").Value = 1.59769778351972
$ws.Range("G13").Value = 0.951822942597957
$ws.Range("H13").Value = 0.9797543412684746
$ws.Range("J13").Value = 0.204103901282533
$ws.Range("K13").Value = 1.554695108138105
$ws.Range("L13").Value = 0.1510483617729115
$ws.Range("M13").Value = 0.5172863027156254
$ws.Range("N13").Value = 1.480956307237449
$ws.Range("O13").Value = 3.907152069197195

$ws.Range("C14").Value = 0.41618424309425
$ws.Range("D14").Value = 0.1839086388174422
$ws.Range("E14").Value = 0.1698614414634427
$ws.Range("F14").Value = 1.597836979162622
$ws.Range("G14").Value = 0.9518935027462589
$ws.Range("H14").Value = 0.9805291516939008
$ws.Range("J14").Value = 0.204288005392808
$ws.Range("K14").Value = 1.533290851272227
$ws.Range("L14").Value = 0.1511495372844891
$ws.Range("M14").Value = 0.5128056769506628
$ws.Range("N14").Value = 1.481536991065568
$ws.Range("O14").Value = 3.908872417453381

$ws.Range("C15").Value = 0.4157001574801882
$ws.Range("D15").Value = 0.183628890501808
$ws.Range("E15").Value = 0.1698448434903383
$ws.Range("F15").Value = 1.597937541312803
$ws.Range("G15").Value = 0.9519476875021411
$ws.Range("H15").Value = 0.9810116849022847
$ws.Range("J15").Value = 0.2044021480975182
$ws.Range("K15").Value = 1.520179848935754
$ws.Range("L15").Value = 0.151212182664402
$ws.Range("M15").Value = 0.5100629071593232
$ws.Range("N15").Value = 1.48190215699239
$ws.Range("O15").Value = 3.909965840848145

$ws.Range("C16").Value = 0.4129776326037131
$ws.Range("D16").Value = 0.1820422488054021
$ws.Range("E16").Value = 0.1697668168530768
$ws.Range("F16").Value = 1.598748416436706
$ws.Range("G16").Value = 0.9524262205775642
$ws.Range("H16").Value = 0.9838986242849757
$ws.Range("J16").Value = 0.2050773204213563
$ws.Range("K16").Value = 1.445043238961148
$ws.Range("L16").Value = 0.151581497521331
$ws.Range("M16").Value = 0.4943722515355233
$ws.Range("N16").Value = 1.484140246380051
$ws.Range("O16").Value = 3.916839711045128

$ws.Range("C17").Value = 0.4113535636573715
$ws.Range("D17").Value = 0.1810836645716449
$ws.Range("E17").Value = 0.1697342252279768
$ws.Range("F17").Value = 1.599455275423367
$ws.Range("G17").Value = 0.9528696751889782
$ws.Range("H17").Value = 0.9857783464976393
$ws.Range("J17").Value = 0.2055103288025926
$ws.Range("K17").Value = 1.398948523282172
$ws.Range("L17").Value = 0.1518172822023676
$ws.Range("M17").Value = 0.4847710443498414
$ws.Range("N17").Value = 1.485643227626127
$ws.Range("O17").Value = 3.921599068011915

$ws.Range("C18").Value = 0.4104364243280259
$ws.Range("D18").Value = 0.180537764102553
$ws.Range("E18").Value = 0.1697211031996595
$ws.Range("F18").Value = 1.59993883264751
$ws.Range("G18").Value = 0.9531798069984063
$ws.Range("H18").Value = 0.986899472018024
$ws.Range("J18").Value = 0.2057663026457739
$ws.Range("K18").Value = 1.372436396386206
$ws.Range("L18").Value = 0.1519562937980083
$ws.Range("M18").Value = 0.4792578331651143
$ws.Range("N18").Value = 1.486555543789294
$ws.Range("O18").Value = 3.924535883277372

$ws.Range("C19").Value = 0.4101288168455994
$ws.Range("D19").Value = 0.1803538702418734
$ws.Range("E19").Value = 0.1697176271724636
$ws.Range("F19").Value = 1.600115783053965
$ws.Range("G19").Value = 0.9532942688664576
$ws.Range("H19").Value = 0.9872859307648554
$ws.Range("J19").Value = 0.2058541598383563
$ws.Range("K19").Value = 1.363459964212154
$ws.Range("L19").Value = 0.1520039442581389
$ws.Range("M19").Value = 0.47739273973324
$ws.Range("N19").Value = 1.486872661985359
$ws.Range("O19").Value = 3.92556447919759

$ws.Range("C20").Value = 0.4115246917957336
$ws.Range("D20").Value = 0.1811851438035177
$ws.Range("E20").Value = 0.169737112849532
$ws.Range("F20").Value = 1.599372061553936
$ws.Range("G20").Value = 0.9528167689555858
$ws.Range("H20").Value = 0.9855741117028032
$ws.Range("J20").Value = 0.2054635183683757
$ws.Range("K20").Value = 1.403855364346896
$ws.Range("L20").Value = 0.1517918313268627
$ws.Range("M20").Value = 0.4857921656014312
$ws.Range("N20").Value = 1.485478282958255
$ws.Range("O20").Value = 3.921071793572338

$ws.Range("C21").Value = 0.416417286321149
$ws.Range("D21").Value = 0.1840430764441976
$ws.Range("E21").Value = 0.1698697025826341
$ws.Range("F21").Value = 1.597792905533396
$ws.Range("G21").Value = 0.9518704924162051
$ws.Range("H21").Value = 0.9802999166264783
$ws.Range("J21").Value = 0.2042336433217997
$ws.Range("K21").Value = 1.539577861644204
$ws.Range("L21").Value = 0.1511196794749416
$ws.Range("M21").Value = 0.5141213829419158
$ws.Range("N21").Value = 1.481364454849597
$ws.Range("O21").Value = 3.908358844313284

$ws.Range("C22").Value = 0.4197653344322703
$ws.Range("D22").Value = 0.1859589335750655
$ws.Range("E22").Value = 0.1700062198372052
$ws.Range("F22").Value = 1.597445749306104
$ws.Range("G22").Value = 0.9517427548083504
$ws.Range("H22").Value = 0.9772083301960208
$ws.Range("J22").Value = 0.2034913437330452
$ws.Range("K22").Value = 1.628269691411106
$ws.Range("L22").Value = 0.1507105015910675
$ws.Range("M22").Value = 0.5327145603477774
$ws.Range("N22").Value = 1.479100338716549
$ws.Range("O22").Value = 3.901825321504845

$ws.Range("C23").Value = 0.4179646887950526
$ws.Range("D23").Value = 0.1849320344600329
$ws.Range("E23").Value = 0.169928806114676
$ws.Range("F23").Value = 1.597568396423071
$ws.Range("G23").Value = 0.951766015015437
$ws.Range("H23").Value = 0.9788258690002039
$ws.Range("J23").Value = 0.2038819018467812
$ws.Range("K23").Value = 1.580934828598117
$ws.Range("L23").Value = 0.1509261378195106
$ws.Range("M23").Value = 0.5227840071784655
$ws.Range("N23").Value = 1.480269969680194
$ws.Range("O23").Value = 3.905149968718575

$ws.Range("C24").Value = 0.4114472731768331
$ws.Range("D24").Value = 0.1811392488325794
$ws.Range("E24").Value = 0.1697357898624716
$ws.Range("F24").Value = 1.599409442142999
$ws.Range("G24").Value = 0.9528405159730511
$ws.Range("H24").Value = 0.9856663202542393
$ws.Range("J24").Value = 0.2054846594616677
$ws.Range("K24").Value = 1.401637017670566
$ws.Range("L24").Value = 0.1518033269020247
$ws.Range("M24").Value = 0.4853304959329279
$ws.Range("N24").Value = 1.4855527041844
$ws.Range("O24").Value = 3.921309549493685

$ws.Range("C25").Value = 0.4050668646405171
$ws.Range("D25").Value = 0.1772568371433323
$ws.Range("E25").Value = 0.1697427248665271
$ws.Range("F25").Value = 1.604339164816153
$ws.Range("G25").Value = 0.9561046766976347
$ws.Range("H25").Value = 0.9945728226670951
$ws.Range("J25").Value = 0.2074780874394762
$ws.Range("K25").Value = 1.208255504161656
$ws.Range("L25").Value = 0.1524994876245795
$ws.Range("M25").Value = 0.4587248717699737
$ws.Range("N25").Value = 1.490299839073856
$ws.Range("O25").Value = 3.936904876889116

